{"js": "// The po_box line in the label template contains a hard-coded literal\n// \"Postfach \" in front of the \"{{ addr.po_box }}\" template tag:\n//   {% endif %}{% if addr.po_box %}Postfach {{ addr.po_box }}\n// That literal word is unnecessary (the commit removes it), so the line\n// should read:\n//   {% endif %}{% if addr.po_box %}{{ addr.po_box }}\nconst body = context.document.body;\n\n// Find the literal run of text \"Postfach \" (including its trailing\n// space) that precedes the \"{{ addr.po_box }}\" Jinja expression.\nconst results = body.search(\"Postfach \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // Remove the matched text entirely (replace with empty string).\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The po_box line in the label template contains a hard-coded literal\n# \"Postfach \" in front of the \"{{ addr.po_box }}\" template tag:\n#   {% endif %}{% if addr.po_box %}Postfach {{ addr.po_box }}\n# That literal word is unnecessary (the commit removes it), so the line\n# should read:\n#   {% endif %}{% if addr.po_box %}{{ addr.po_box }}\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Postfach \"\n$find.Replacement.Text = \"\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdFindContinue (1) wrap mode, wdReplaceAll (2) replace mode.\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
